$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be stored as text so values like "7.30" or
# "0.0000158" keep their exact original formatting instead of being
# auto-converted into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '56.809.33'
$ws.Range("E2").Value = '  -0.88%  '

$ws.Range("D3").Value = '2.967.13'
$ws.Range("E3").Value = '  -1.80%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '497.42'
$ws.Range("E5").Value = '  -4.04%  '

$ws.Range("D6").Value = '136.87'
$ws.Range("E6").Value = '  -3.45%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("E8").Value = '  -2.39%  '

$ws.Range("D9").Value = '7.30'
$ws.Range("E9").Value = '  -3.73%  '

$ws.Range("E10").Value = '  -2.43%  '

$ws.Range("D11").Value = '0.356'
$ws.Range("E11").Value = '  -1.20%  '

$ws.Range("D12").Value = '3.476.01'
$ws.Range("E12").Value = '  -1.92%  '

$ws.Range("E13").Value = '  -2.00%  '

$ws.Range("D14").Value = '25.86'
$ws.Range("E14").Value = '  -0.99%  '

$ws.Range("D15").Value = '0.0000158'
$ws.Range("E15").Value = '  -1.83%  '

$ws.Range("D16").Value = '56.888.81'
$ws.Range("E16").Value = '  -0.76%  '

$ws.Range("D17").Value = '6.03'
$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("D18").Value = '2.967.83'
$ws.Range("E18").Value = '  -1.93%  '

$ws.Range("D19").Value = '12.51'
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").Value = '7.80'
$ws.Range("E20").Value = '  -2.07%  '

$ws.Range("D21").Value = '317.61'
$ws.Range("E21").Value = '  -4.15%  '

$ws.Range("E22").Value = '  -0.12%  '

$ws.Range("E23").Value = '  -0.82%  '

$ws.Range("D24").Value = '0.484'
$ws.Range("E24").Value = '  -0.87%  '

$ws.Range("D25").Value = '63.31'
$ws.Range("E25").Value = '  -1.45%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").Value = '0.163'
$ws.Range("E27").Value = '  -5.47%  '

$ws.Range("D28").Value = '0.0₃0886'
$ws.Range("E28").Value = '  -4.27%  '

$ws.Range("D29").Value = '6.47'
$ws.Range("E29").Value = '  -4.65%  '

$ws.Range("D30").Value = '7.00'
$ws.Range("E30").Value = '  -2.96%  '

$ws.Range("E31").Value = '  -3.91%  '

$ws.Range("E32").Value = '  -6.59%  '

$ws.Range("D33").Value = '20.07'
$ws.Range("E33").Value = '  -3.49%  '

$ws.Range("D34").Value = '156.52'
$ws.Range("E34").Value = '  -1.38%  '

$ws.Range("D35").Value = '4.55'
$ws.Range("E35").Value = '  -2.26%  '

$ws.Range("D36").Value = '5.72'
$ws.Range("E36").Value = '  -1.15%  '

$ws.Range("E37").Value = '  -5.08%  '

$ws.Range("D38").Value = '23.94'
$ws.Range("E38").Value = '  -2.35%  '

$ws.Range("D39").Value = '0.0662'
$ws.Range("E39").Value = '  -2.34%  '

$ws.Range("D40").Value = '2.998.36'
$ws.Range("E40").Value = '  -1.96%  '

$ws.Range("D41").Value = '37.51'
$ws.Range("E41").Value = '  +0.17%  '

$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("D43").Value = '3.71'
$ws.Range("E43").Value = '  -0.97%  '

$ws.Range("D44").Value = '0.636'
$ws.Range("E44").Value = '  -2.83%  '

$ws.Range("D45").Value = '2.191.11'
$ws.Range("E45").Value = '  -4.89%  '

$ws.Range("E46").Value = '  -4.70%  '

$ws.Range("D47").Value = '5.91'
$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("D48").Value = '0.933'
$ws.Range("E48").Value = '  -8.04%  '

$ws.Range("E49").Value = '  -3.83%  '

$ws.Range("D50").Value = '19.09'
$ws.Range("E50").Value = '  -2.25%  '

$ws.Range("D51").Value = '1.79'
$ws.Range("E51").Value = '  -11.74%  '
